$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 347
$ws.Range("I92").Value = 356.7143
$ws.Range("J92").Value = 301.66666
$ws.Range("K92").Value = 356.7143
$ws.Range("L92").Value = 301.66666
$ws.Range("M92").Value = 891.2857
$ws.Range("N92").Value = -2797.66666
$ws.Range("H129").Value = 751.5208
$ws.Range("J129").Value = 799.67444
$ws.Range("L129").Value = 2399.02332
$ws.Range("N129").Value = -12399.02332
$ws.Range("H138").Value = 2014.5466
$ws.Range("I138").Value = 1110.6
$ws.Range("J138").Value = 2617.1777
$ws.Range("K138").Value = 3331.8
$ws.Range("L138").Value = 7851.533100000001
$ws.Range("M138").Value = 1808.2
$ws.Range("N138").Value = -18131.5331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1706.5686
$ws.Range("I61").Value = 1517.8182
$ws.Range("J61").Value = 2052.611
$ws.Range("K61").Value = 1517.8182
$ws.Range("L61").Value = 2052.611
$ws.Range("M61").Value = -1305.8182
$ws.Range("N61").Value = -2476.611
$ws.Range("H119").Value = 15600
$ws.Range("J119").Value = 15600
$ws.Range("L119").Value = 15600
$ws.Range("N119").Value = -25276
$ws.Range("H124").Value = 11809.333
$ws.Range("J124").Value = 11809.333
$ws.Range("L124").Value = 11809.333
$ws.Range("N124").Value = -21629.333
$ws.Range("H125").Value = 31428.666
$ws.Range("J125").Value = 31428.666
$ws.Range("L125").Value = 31428.666
$ws.Range("N125").Value = -41268.666
$ws.Range("H136").Value = 1706.5686
$ws.Range("I136").Value = 1517.8182
$ws.Range("J136").Value = 2052.611
$ws.Range("K136").Value = 4553.4546
$ws.Range("L136").Value = 6157.833
$ws.Range("M136").Value = -2003.4546
$ws.Range("N136").Value = -11257.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 952.94116
$ws.Range("I94").Value = 864.4545000000001
$ws.Range("J94").Value = 1115.1666
$ws.Range("K94").Value = 864.4545000000001
$ws.Range("L94").Value = 1115.1666
$ws.Range("M94").Value = -413.4545000000001
$ws.Range("N94").Value = -2017.1666
$ws.Range("H110").Value = 31700
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 31700
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 31700
$ws.Range("M110").Value = $null
$ws.Range("N110").Value = -39880
$ws.Range("H130").Value = 43284
$ws.Range("J130").Value = 43284
$ws.Range("L130").Value = 43284
$ws.Range("N130").Value = -53324

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4071.375
$ws.Range("I31").Value = 1075.9166
$ws.Range("K31").Value = 1075.9166
$ws.Range("M31").Value = -780.9166
$ws.Range("H34").Value = 4071.375
$ws.Range("I34").Value = 1075.9166
$ws.Range("K34").Value = 1075.9166
$ws.Range("M34").Value = -873.9166
$ws.Range("H99").Value = 3639.0435
$ws.Range("I99").Value = 2655.4443
$ws.Range("K99").Value = 2655.4443
$ws.Range("M99").Value = -1157.4443
$ws.Range("H107").Value = 1108.4642
$ws.Range("I107").Value = 443.06668
$ws.Range("J107").Value = 1876.2307
$ws.Range("K107").Value = 443.06668
$ws.Range("L107").Value = 1876.2307
$ws.Range("M107").Value = 1476.93332
$ws.Range("N107").Value = -5716.2307
$ws.Range("H124").Value = 10151.8
$ws.Range("I124").Value = 8948.357
$ws.Range("J124").Value = 27000
$ws.Range("K124").Value = 8948.357
$ws.Range("L124").Value = 27000
$ws.Range("M124").Value = -6493.357
$ws.Range("N124").Value = -31910
$ws.Range("H126").Value = 3639.0435
$ws.Range("I126").Value = 2655.4443
$ws.Range("K126").Value = 7966.3329
$ws.Range("M126").Value = -5496.3329
$ws.Range("H137").Value = 26936.924
$ws.Range("I137").Value = 9000
$ws.Range("K137").Value = 9000
$ws.Range("M137").Value = -3900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 20833776
$ws.Range("I92").Value = 41666830
$ws.Range("K92").Value = 125000490
$ws.Range("M92").Value = -124999242
$ws.Range("H107").Value = 7006.143
$ws.Range("J107").Value = 304.44446
$ws.Range("L107").Value = 913.33338
$ws.Range("N107").Value = -4753.33338
$ws.Range("H131").Value = 756.17
$ws.Range("J131").Value = 768.3402
$ws.Range("L131").Value = 2305.0206
$ws.Range("N131").Value = -12385.0206

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5109.95
$ws.Range("I7").Value = 4979.933
$ws.Range("J7").Value = 5500
$ws.Range("K7").Value = 4979.933
$ws.Range("L7").Value = 5500
$ws.Range("M7").Value = -4867.933
$ws.Range("N7").Value = -5724
$ws.Range("H61").Value = 4141.1665
$ws.Range("I61").Value = 2538.9
$ws.Range("J61").Value = 12152.5
$ws.Range("K61").Value = 2538.9
$ws.Range("L61").Value = 12152.5
$ws.Range("M61").Value = -2336.9
$ws.Range("N61").Value = -12556.5
$ws.Range("H93").Value = 2039.1111
$ws.Range("I93").Value = 1846.9333
$ws.Range("K93").Value = 1846.9333
$ws.Range("M93").Value = -598.9332999999999
$ws.Range("H113").Value = 4141.1665
$ws.Range("I113").Value = 2538.9
$ws.Range("J113").Value = 12152.5
$ws.Range("K113").Value = 2538.9
$ws.Range("L113").Value = 12152.5
$ws.Range("M113").Value = -368.9000000000001
$ws.Range("N113").Value = -16492.5
$ws.Range("H122").Value = 894401.25
$ws.Range("I122").Value = 1156360.5
$ws.Range("K122").Value = 3469081.5
$ws.Range("M122").Value = -3466631.5
$ws.Range("H126").Value = 5109.95
$ws.Range("I126").Value = 4979.933
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 14939.799
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -12469.799
$ws.Range("N126").Value = -21440
$ws.Range("H132").Value = 637478.3
$ws.Range("I132").Value = 805339.3
$ws.Range("K132").Value = 2416017.9
$ws.Range("M132").Value = -2413487.9
$ws.Range("H136").Value = 1952.8148
$ws.Range("I136").Value = 1859.24
$ws.Range("K136").Value = 5577.72
$ws.Range("M136").Value = -3027.72

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 83333670
$ws.Range("I107").Value = 111111350
$ws.Range("J107").Value = 636.6667
$ws.Range("K107").Value = 333334050
$ws.Range("L107").Value = 1910.0001
$ws.Range("M107").Value = -333332130
$ws.Range("N107").Value = -5750.0001
$ws.Range("H132").Value = 1351.909
$ws.Range("I132").Value = 764.44446
$ws.Range("J132").Value = 2056.8667
$ws.Range("K132").Value = 2293.33338
$ws.Range("L132").Value = 6170.6001
$ws.Range("M132").Value = 236.66662
$ws.Range("N132").Value = -11230.6001
$ws.Range("H136").Value = 27167594
$ws.Range("I136").Value = 34409890
$ws.Range("K136").Value = 103229670
$ws.Range("M136").Value = -103227120
